# Update the weekly duty-roster tables on both sheets to match the
# re-uploaded schedule. Names are addressed explicitly by sheet + cell so
# the edit is independent of whichever sheet happens to be "active".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("schedule_2021_7")
$ws2 = $wb.Worksheets.Item("backup_2021_7")

# ---- schedule_2021_7 ----
$ws1.Range("B2").Value  = "Barna Dávid"
$ws1.Range("D2").Value  = "Nemes Péter"
$ws1.Range("E2").Value  = "Kovács Gusztáv Márk"
$ws1.Range("F2").Value  = "Gregor Patrik"
$ws1.Range("H2").Value  = "Gregor Patrik"

$ws1.Range("B3").Value  = "Földi Bence"
$ws1.Range("D3").Value  = "Ördög Márk"
$ws1.Range("E3").Value  = "Hugyecz Árpád"
$ws1.Range("F3").Value  = "Kovács Gusztáv Márk"
$ws1.Range("G3").Value  = "Nagy Roland"
$ws1.Range("H3").Value  = "Nagy Roland"

$ws1.Range("D4").Value  = "Hugyecz Árpád"
$ws1.Range("E4").Value  = "Biró Levente"
$ws1.Range("G4").Value  = "Ördög Márk"
$ws1.Range("H4").Value  = "Ördög Márk"

$ws1.Range("B7").Value  = "Czimmer Sándor"
$ws1.Range("C7").Value  = "Czimmer Sándor"
$ws1.Range("D7").Value  = "Koleszár Zoltán"
$ws1.Range("E7").Value  = "Czimmer Sándor"
$ws1.Range("G7").Value  = "Balla Mihály"

$ws1.Range("B8").Value  = "Nagy Dávid"
$ws1.Range("C8").Value  = "Hajdu Krisztián"
$ws1.Range("D8").Value  = "Kormányos Patrik"
$ws1.Range("E8").Value  = "Gregor Patrik"
$ws1.Range("F8").Value  = "Hajdu Krisztián"
$ws1.Range("G8").Value  = "Barna Dávid"
$ws1.Range("H8").Value  = "Hajdu Krisztián"

$ws1.Range("E9").Value  = "Paragi Gábor"
$ws1.Range("F9").Value  = "Kormányos Patrik"
$ws1.Range("G9").Value  = "Gregor Patrik"
$ws1.Range("H9").Value  = "Kormányos Patrik"

$ws1.Range("B10").Value = "Emődi Máté"
$ws1.Range("C10").Value = "Peti Márk"
$ws1.Range("D10").Value = "Peti Márk"
$ws1.Range("E10").Value = "Peti Márk"
$ws1.Range("F10").Value = "Paragi Gábor"
$ws1.Range("G10").Value = "Hajdu Krisztián"
$ws1.Range("H10").Value = "Paragi Gábor"

# ---- backup_2021_7 ----
$ws2.Range("B2").Value = "Biró Levente"
$ws2.Range("H4").Value = "Biró Levente"
$ws2.Range("C8").Value = "Emődi Máté"
$ws2.Range("D9").Value = "Emődi Máté"
